$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and the date range banner) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Crime statistics table updates (rows 14-29) ---

# A handful of cells change data TYPE in this edit (numeric <-> text
# placeholder). Handle those first so their resulting style matches the
# convention already used elsewhere in the sheet: style 14/General is the
# "no data" text placeholder (e.g. C14), style 15/16 hold real numbers.

# numeric -> text "0" placeholder cells
$ws.Range("C15").Formula = "'0"
$ws.Range("C26").Formula = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# text placeholder -> numeric cells (restore the plain numeric format this
# column uses elsewhere before writing the number)
$ws.Range("D22").NumberFormat = $ws.Range("C16").NumberFormat()
$ws.Range("E22").NumberFormat = $ws.Range("E16").NumberFormat()
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300

# Row 14
$ws.Range("M14").Value = 50

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -25
$ws.Range("J15").Value = 26
$ws.Range("K15").Value = 34.615384615384
$ws.Range("L15").Value = 45.833333333333

# Row 16
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 57.142857142857
$ws.Range("F16").Value = 32
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 88.235294117647
$ws.Range("I16").Value = 261
$ws.Range("J16").Value = 207
$ws.Range("K16").Value = 26.086956521739
$ws.Range("L16").Value = 47.457627118644
$ws.Range("M16").Value = -13
$ws.Range("N16").Value = -76.633840644583

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 11.538461538461
$ws.Range("I17").Value = 361
$ws.Range("J17").Value = 329
$ws.Range("K17").Value = 9.726443768996
$ws.Range("L17").Value = 36.742424242424
$ws.Range("M17").Value = 28.014184397163
$ws.Range("N17").Value = -2.168021680216

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 128
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = 7.563025210084
$ws.Range("L18").Value = -12.328767123287
$ws.Range("M18").Value = -49.803921568627
$ws.Range("N18").Value = -92.426035502958

# Row 19
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 108.333333333333
$ws.Range("F19").Value = 91
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 59.649122807017
$ws.Range("I19").Value = 841
$ws.Range("J19").Value = 466
$ws.Range("K19").Value = 80.472103004291
$ws.Range("L19").Value = 86.474501108647
$ws.Range("M19").Value = 98.349056603773
$ws.Range("N19").Value = -31.902834008097

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 24
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 41.176470588235
$ws.Range("I20").Value = 267
$ws.Range("J20").Value = 179
$ws.Range("K20").Value = 49.162011173184
$ws.Range("L20").Value = 53.448275862069
$ws.Range("M20").Value = 39.79057591623
$ws.Range("N20").Value = -86.20867768595

# Row 21
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 54.545454545454
$ws.Range("F21").Value = 189
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = 34.042553191489
$ws.Range("I21").Value = 1896
$ws.Range("J21").Value = 1330
$ws.Range("K21").Value = 42.556390977443
$ws.Range("L21").Value = 52.044907778668
$ws.Range("M21").Value = 27.935222672064
$ws.Range("N21").Value = -70.323994365315

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 58
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = 132
$ws.Range("L22").Value = 262.5
$ws.Range("M22").Value = 114.814814814815

# Row 24
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = 21.428571428571
$ws.Range("F24").Value = 188
$ws.Range("G24").Value = 125
$ws.Range("H24").Value = 50.4
$ws.Range("I24").Value = 1629
$ws.Range("J24").Value = 1169
$ws.Range("K24").Value = 39.349871685201
$ws.Range("L24").Value = 34.962717481358
$ws.Range("M24").Value = 67.248459958932

# Row 25
$ws.Range("D25").Value = 20
$ws.Range("E25").Value = -15
$ws.Range("F25").Value = 79
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = -8.13953488372
$ws.Range("I25").Value = 771
$ws.Range("J25").Value = 725
$ws.Range("K25").Value = 6.344827586206
$ws.Range("L25").Value = 19.53488372093
$ws.Range("M25").Value = -4.342431761786

# Row 26
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = -33.333333333333
$ws.Range("J26").Value = 46
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 17.948717948717

# Row 27
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 800
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 18.181818181818
$ws.Range("I27").Value = 94
$ws.Range("J27").Value = 92
$ws.Range("K27").Value = 2.173913043478
$ws.Range("L27").Value = 28.767123287671

# Row 28
$ws.Range("M28").Value = 0

# Row 29
$ws.Range("M29").Value = 16.666666666666

